$d = $word.ActiveDocument

# --- Step 1: Append sentence to the "esta diferencia:" paragraph ---
$p = $d.Paragraphs(4)
$r = $p.Range
$r.End = $r.End - 1  # exclude paragraph mark
$r.InsertAfter(" La membresía Estándar tiene un volumen de cerca de 40. En comparación, las membresías de Premium y Exclusive combinados apenas llegan a un volumen de 30.")

Write-Output "Step1 done"

# --- Step 2: Replace paragraph 5 content entirely (was "La membresía Estándar..." list of runs) ---
$p = $d.Paragraphs(5)
$r = $p.Range
$r.End = $r.End - 1  # exclude paragraph mark
$r.Text = ""
$r.Collapse(1)  # collapse to start

$r.InsertAfter("La razón principal de que la cuenta Estándar sea la más utilizada es su facilidad de acceso. Este nivel suele ser la opción gratuita o la más económica, lo que elimina el obstáculo inicial del pago para los nuevos usuarios. Al ofrecer las funciones básicas necesarias, el nivel Estándar se convierte en el camino inicial y suficiente para la mayoría, evitando que los usuarios no encuentren razones suficientes para mejorar su membresía. ")
$r.Collapse(0)

$boldStart = $r.Start
$r.InsertAfter("Sin embargo")
$boldEnd = $r.End
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Bold = 1
$r.Collapse(0)

$r.InsertAfter(", es importante destacar que la suma de usuarios Premium y Exclusive representa casi un 30% del total, lo que indica una tasa de conversión saludable. Esto sugiere que, aunque la entrada es masiva en Estándar, existe un segmento fiel dispuesto a pagar por eventos exclusivos o acceso total, validando nuestro modelo de negocio que consiste en atraer a las masas para captar mas usuarios que paguen.")
$r.Collapse(0)

Write-Output "Step2 done"

# --- Step 3: Clear paragraph 6 (old "La razón principal..." duplicate) to empty ---
$p = $d.Paragraphs(6)
$r = $p.Range
$r.End = $r.End - 1  # exclude paragraph mark
$r.Text = ""

Write-Output "Step3 done"

# --- Step 4: Add left indent (1080 twips = 54pt) to "La actividad del museo..." paragraph ---
$p = $d.Paragraphs(8)
$p.Range.ParagraphFormat.LeftIndent = 54

Write-Output "Step4 done"

# --- Step 5: Split the "El gráfico de líneas..." paragraph after "cerca de cero." ---
$p = $d.Paragraphs(9)
$r = $p.Range.Duplicate
$null = $r.Find.Execute("cerca de cero.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($r.End, $r.End)
$splitPoint.InsertParagraphAfter()

Write-Output "Step5a done"

# --- Step 5b: Replace paragraph 10's content (old continuation) with new text ---
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range
$r10.End = $r10.End - 1  # exclude paragraph mark, keep leading space before "Este" removed too since we'll set full text
$r10.Text = ""
$r10.Collapse(1)

$r10.InsertAfter("Este comportamiento y las modificaciones que se generan en la línea a lo largo del tiempo indican que la plataforma depende significativamente de eventos puntuales. Esto confirma que nuestra ")
$r10.Collapse(0)

$underlineStart = $r10.Start
$r10.InsertAfter("audi")
$underlineEnd1 = $r10.End
$uRange1 = $d.Range($underlineStart, $underlineEnd1)
$uRange1.Font.Underline = 1
$r10.Collapse(0)

$bmRange = $d.Range($r10.Start, $r10.Start)
$bmRange.Bookmarks.Add("_GoBack")

$underline2Start = $r10.Start
$r10.InsertAfter("encia")
$underline2End = $r10.End
$uRange2 = $d.Range($underline2Start, $underline2End)
$uRange2.Font.Underline = 1
$r10.Collapse(0)

$r10.InsertAfter(" es ""reactiva"": ingresan masivamente ante un estímulo pero no tienen un hábito de visita orgánica diaria. Esto nos muestra que habría que buscar la forma de lograr que las visitas al museo se den de una manera más regular y no dependan exclusivamente de sus eventos, quizás implementando dinámicas de retención diaria para evitar los ""días muertos"" que se observan entre pico y pico")
$r10.Collapse(0)
$r10.InsertAfter(".")
$r10.Collapse(0)

Write-Output "Step5b done"
